$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RawData")

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "First"
$ws.Range("C1").Value = "Last"
$ws.Range("D1").Value = "School"
$ws.Range("E1").Value = "SMCS"
$ws.Range("F1").Value = "Global"
$ws.Range("G1").Value = "Humanities"

$ws.Activate()
$ws.Range("F2").Select()
